$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 3, pushing the existing data (rows 3..37) down
# to rows 4..38. This mirrors the weekly refresh: a brand new record lands
# at row 3 and everything else shifts down by one row.
$ws.Rows.Item(3).Insert("xlShiftDown")

# Populate the newly inserted row 3 with this week's new record.
$ws.Cells.Item(3, 1).Value = 7
$ws.Cells.Item(3, 2).Value = 'Terminal Hortofrutícola Agro Chillán'
$ws.Cells.Item(3, 3).Value = 'Ñuble'
$ws.Cells.Item(3, 4).Value = 44761
$ws.Cells.Item(3, 5).Value = 16
$ws.Cells.Item(3, 6).Value = 100112013
$ws.Cells.Item(3, 7).Value = 'Alcachofa'
$ws.Cells.Item(3, 8).Value = 'Argentina(o)'
$ws.Cells.Item(3, 9).Value = 'Primera'
$ws.Cells.Item(3, 10).Value = 60
$ws.Cells.Item(3, 11).Value = 15000
$ws.Cells.Item(3, 12).Value = 16000
$ws.Cells.Item(3, 13).Value = 15500
$ws.Cells.Item(3, 14).Value = '$/caja 50 unidades'
$ws.Cells.Item(3, 15).Value = 'Provincia de Limarí'
$ws.Cells.Item(3, 16).Value = 310
$ws.Cells.Item(3, 17).Value = 50
$ws.Cells.Item(3, 18).Value = 'Hortaliza'
